$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text storage for the Price column (D) so numeric-looking
# strings such as "0.999" or "51.294.98" are not auto-converted to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.294.98'
$ws.Range("E2").Value = '  -1.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.971.98'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '382.35'
$ws.Range("E5").Value = '  +1.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.99'
$ws.Range("E6").Value = '  -4.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.540'
$ws.Range("E7").Value = '  -1.15%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.588'
$ws.Range("E9").Value = '  -1.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.82'
$ws.Range("E10").Value = '  -2.29%  '
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0841'
$ws.Range("E12").Value = '  -0.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.441.00'
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.14'
$ws.Range("E14").Value = '  -2.33%  '
$ws.Range("E15").Value = '  +0.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.969.82'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.992'
$ws.Range("E17").Value = '  +4.45%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.202.46'
$ws.Range("E18").Value = '  -1.12%  '
$ws.Range("E19").Value = '  -3.19%  '
$ws.Range("E20").Value = '  -0.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.76'
$ws.Range("E21").Value = '  -3.24%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.80'
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '260.51'
$ws.Range("E24").Value = '  -1.04%  '
$ws.Range("E25").Value = '  +5.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.09'
$ws.Range("E26").Value = '  +11.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.40'
$ws.Range("E27").Value = '  +6.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.116'
$ws.Range("E28").Value = '  +10.84%  '
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.168'
$ws.Range("E29").Value = '  -2.20%  '
$ws.Range("B30").Value = 'LEO'
$ws.Range("C30").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.13'
$ws.Range("E30").Value = '  -0.48%  '
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("E32").Value = '  -1.10%  '
$ws.Range("E33").Value = '  -1.69%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '34.28'
$ws.Range("E34").Value = '  -2.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.79'
$ws.Range("E35").Value = '  +0.53%  '
$ws.Range("E36").Value = '  -2.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0450'
$ws.Range("E37").Value = '  +3.69%  '
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("E39").Value = '  -3.66%  '
$ws.Range("E40").Value = '  -2.39%  '
$ws.Range("E41").Value = '  -3.70%  '
$ws.Range("E43").Value = '  -3.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '122.69'
$ws.Range("E44").Value = '  +1.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.44'
$ws.Range("E45").Value = '  -4.55%  '
$ws.Range("E46").Value = '  -1.09%  '
$ws.Range("E47").Value = '  -2.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.36'
$ws.Range("E48").Value = '  +2.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.027.78'
$ws.Range("E49").Value = '  -1.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.25'
$ws.Range("E50").Value = '  +0.38%  '
$ws.Range("E51").Value = '  -1.89%  '
